$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI TPM numbers were regenerated, the target-cluster list grew
# an "ECs" category and lost "Neutrophils", and "Resolving-Mac" is no longer a
# sending cluster - so the table shrinks from 12 data rows to 10.
# Drop the two now-unused trailing data rows first (old rows 12 and 13).
$ws.Range("A12:A13").EntireRow.Delete()

# Refresh rows 2-11 (FAPs/Inflammatory-Mac senders x ECs/FAPs/Inflammatory-Mac/MuSCs/
# Resolving-Mac targets) with the new TPM-derived values, column by column.
# Column A
$ws.Cells.Item(2,1).Value2 = "FAPs"
$ws.Cells.Item(3,1).Value2 = "FAPs"
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(7,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(8,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(9,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(10,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(11,1).Value2 = "Inflammatory-Mac"

# Column B
$ws.Cells.Item(2,2).Value2 = "Rspo1"
$ws.Cells.Item(3,2).Value2 = "Rspo1"
$ws.Cells.Item(4,2).Value2 = "Rspo1"
$ws.Cells.Item(5,2).Value2 = "Rspo1"
$ws.Cells.Item(6,2).Value2 = "Rspo1"
$ws.Cells.Item(7,2).Value2 = "Rspo1"
$ws.Cells.Item(8,2).Value2 = "Rspo1"
$ws.Cells.Item(9,2).Value2 = "Rspo1"
$ws.Cells.Item(10,2).Value2 = "Rspo1"
$ws.Cells.Item(11,2).Value2 = "Rspo1"

# Column C
$ws.Cells.Item(2,3).Value2 = "Lgr5"
$ws.Cells.Item(3,3).Value2 = "Lgr5"
$ws.Cells.Item(4,3).Value2 = "Lgr5"
$ws.Cells.Item(5,3).Value2 = "Lgr5"
$ws.Cells.Item(6,3).Value2 = "Lgr5"
$ws.Cells.Item(7,3).Value2 = "Lgr5"
$ws.Cells.Item(8,3).Value2 = "Lgr5"
$ws.Cells.Item(9,3).Value2 = "Lgr5"
$ws.Cells.Item(10,3).Value2 = "Lgr5"
$ws.Cells.Item(11,3).Value2 = "Lgr5"

# Column D
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(4,4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(5,4).Value2 = "MuSCs"
$ws.Cells.Item(6,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(7,4).Value2 = "ECs"
$ws.Cells.Item(8,4).Value2 = "FAPs"
$ws.Cells.Item(9,4).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(10,4).Value2 = "MuSCs"
$ws.Cells.Item(11,4).Value2 = "Resolving-Mac"

# Column E
$ws.Cells.Item(2,5).Value2 = 2.0
$ws.Cells.Item(3,5).Value2 = 2.0
$ws.Cells.Item(4,5).Value2 = 2.0
$ws.Cells.Item(5,5).Value2 = 2.0
$ws.Cells.Item(6,5).Value2 = 2.0
$ws.Cells.Item(7,5).Value2 = 3.0
$ws.Cells.Item(8,5).Value2 = 3.0
$ws.Cells.Item(9,5).Value2 = 3.0
$ws.Cells.Item(10,5).Value2 = 3.0
$ws.Cells.Item(11,5).Value2 = 3.0

# Column F
$ws.Cells.Item(2,6).Value2 = 0.6666666666666666
$ws.Cells.Item(3,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(5,6).Value2 = 0.6666666666666666
$ws.Cells.Item(6,6).Value2 = 0.6666666666666666
$ws.Cells.Item(7,6).Value2 = 1.0
$ws.Cells.Item(8,6).Value2 = 1.0
$ws.Cells.Item(9,6).Value2 = 1.0
$ws.Cells.Item(10,6).Value2 = 1.0
$ws.Cells.Item(11,6).Value2 = 1.0

# Column G
$ws.Cells.Item(2,7).Value2 = 0.07032533333333334
$ws.Cells.Item(3,7).Value2 = 0.07032533333333334
$ws.Cells.Item(4,7).Value2 = 0.07032533333333334
$ws.Cells.Item(5,7).Value2 = 0.07032533333333334
$ws.Cells.Item(6,7).Value2 = 0.07032533333333334
$ws.Cells.Item(7,7).Value2 = 0.3024513333333333
$ws.Cells.Item(8,7).Value2 = 0.3024513333333333
$ws.Cells.Item(9,7).Value2 = 0.3024513333333333
$ws.Cells.Item(10,7).Value2 = 0.3024513333333333
$ws.Cells.Item(11,7).Value2 = 0.3024513333333333

# Column H
$ws.Cells.Item(2,8).Value2 = 0.210976
$ws.Cells.Item(3,8).Value2 = 0.210976
$ws.Cells.Item(4,8).Value2 = 0.210976
$ws.Cells.Item(5,8).Value2 = 0.210976
$ws.Cells.Item(6,8).Value2 = 0.210976
$ws.Cells.Item(7,8).Value2 = 0.907354
$ws.Cells.Item(8,8).Value2 = 0.907354
$ws.Cells.Item(9,8).Value2 = 0.907354
$ws.Cells.Item(10,8).Value2 = 0.907354
$ws.Cells.Item(11,8).Value2 = 0.907354

# Column I
$ws.Cells.Item(2,9).Value2 = 0.1886527232569993
$ws.Cells.Item(3,9).Value2 = 0.1886527232569993
$ws.Cells.Item(4,9).Value2 = 0.1886527232569993
$ws.Cells.Item(5,9).Value2 = 0.1886527232569993
$ws.Cells.Item(6,9).Value2 = 0.1886527232569993
$ws.Cells.Item(7,9).Value2 = 0.8113472767430007
$ws.Cells.Item(8,9).Value2 = 0.8113472767430007
$ws.Cells.Item(9,9).Value2 = 0.8113472767430007
$ws.Cells.Item(10,9).Value2 = 0.8113472767430007
$ws.Cells.Item(11,9).Value2 = 0.8113472767430007

# Column J
$ws.Cells.Item(2,10).Value2 = 0.1886527232569993
$ws.Cells.Item(3,10).Value2 = 0.1886527232569993
$ws.Cells.Item(4,10).Value2 = 0.1886527232569993
$ws.Cells.Item(5,10).Value2 = 0.1886527232569993
$ws.Cells.Item(6,10).Value2 = 0.1886527232569993
$ws.Cells.Item(7,10).Value2 = 0.8113472767430007
$ws.Cells.Item(8,10).Value2 = 0.8113472767430007
$ws.Cells.Item(9,10).Value2 = 0.8113472767430007
$ws.Cells.Item(10,10).Value2 = 0.8113472767430007
$ws.Cells.Item(11,10).Value2 = 0.8113472767430007

# Column K
$ws.Cells.Item(2,11).Value2 = 2.0
$ws.Cells.Item(3,11).Value2 = 3.0
$ws.Cells.Item(4,11).Value2 = 1.0
$ws.Cells.Item(5,11).Value2 = 3.0
$ws.Cells.Item(6,11).Value2 = 1.0
$ws.Cells.Item(7,11).Value2 = 2.0
$ws.Cells.Item(8,11).Value2 = 3.0
$ws.Cells.Item(9,11).Value2 = 1.0
$ws.Cells.Item(10,11).Value2 = 3.0
$ws.Cells.Item(11,11).Value2 = 1.0

# Column L
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,12).Value2 = 1.0
$ws.Cells.Item(4,12).Value2 = 0.3333333333333333
$ws.Cells.Item(5,12).Value2 = 1.0
$ws.Cells.Item(6,12).Value2 = 0.3333333333333333
$ws.Cells.Item(7,12).Value2 = 0.6666666666666666
$ws.Cells.Item(8,12).Value2 = 1.0
$ws.Cells.Item(9,12).Value2 = 0.3333333333333333
$ws.Cells.Item(10,12).Value2 = 1.0
$ws.Cells.Item(11,12).Value2 = 0.3333333333333333

# Column M
$ws.Cells.Item(2,13).Value2 = 0.184159
$ws.Cells.Item(3,13).Value2 = 0.6233136666666667
$ws.Cells.Item(4,13).Value2 = 0.064319
$ws.Cells.Item(5,13).Value2 = 0.4910533333333333
$ws.Cells.Item(6,13).Value2 = 0.02060433333333333
$ws.Cells.Item(7,13).Value2 = 0.184159
$ws.Cells.Item(8,13).Value2 = 0.6233136666666667
$ws.Cells.Item(9,13).Value2 = 0.064319
$ws.Cells.Item(10,13).Value2 = 0.4910533333333333
$ws.Cells.Item(11,13).Value2 = 0.02060433333333333

# Column N
$ws.Cells.Item(2,14).Value2 = 0.5524770000000001
$ws.Cells.Item(3,14).Value2 = 1.869941
$ws.Cells.Item(4,14).Value2 = 0.192957
$ws.Cells.Item(5,14).Value2 = 1.47316
$ws.Cells.Item(6,14).Value2 = 0.061813
$ws.Cells.Item(7,14).Value2 = 0.5524770000000001
$ws.Cells.Item(8,14).Value2 = 1.869941
$ws.Cells.Item(9,14).Value2 = 0.192957
$ws.Cells.Item(10,14).Value2 = 1.47316
$ws.Cells.Item(11,14).Value2 = 0.061813

# Column O
$ws.Cells.Item(2,15).Value2 = 0.1331158254681294
$ws.Cells.Item(3,15).Value2 = 0.4505504116763221
$ws.Cells.Item(4,15).Value2 = 0.04649176406412185
$ws.Cells.Item(5,15).Value2 = 0.3549485488927676
$ws.Cells.Item(6,15).Value2 = 0.0148934498986591
$ws.Cells.Item(7,15).Value2 = 0.1331158254681294
$ws.Cells.Item(8,15).Value2 = 0.4505504116763221
$ws.Cells.Item(9,15).Value2 = 0.04649176406412185
$ws.Cells.Item(10,15).Value2 = 0.3549485488927676
$ws.Cells.Item(11,15).Value2 = 0.0148934498986591

# Column P
$ws.Cells.Item(2,16).Value2 = 0.1331158254681294
$ws.Cells.Item(3,16).Value2 = 0.450550411676322
$ws.Cells.Item(4,16).Value2 = 0.04649176406412184
$ws.Cells.Item(5,16).Value2 = 0.3549485488927676
$ws.Cells.Item(6,16).Value2 = 0.0148934498986591
$ws.Cells.Item(7,16).Value2 = 0.1331158254681294
$ws.Cells.Item(8,16).Value2 = 0.450550411676322
$ws.Cells.Item(9,16).Value2 = 0.04649176406412184
$ws.Cells.Item(10,16).Value2 = 0.3549485488927676
$ws.Cells.Item(11,16).Value2 = 0.0148934498986591

# Column Q
$ws.Cells.Item(2,17).Value2 = 0.01295104306133334
$ws.Cells.Item(3,17).Value2 = 0.04383474137955556
$ws.Cells.Item(4,17).Value2 = 0.004523255114666667
$ws.Cells.Item(5,17).Value2 = 0.03453348935111111
$ws.Cells.Item(6,17).Value2 = 0.001449006609777778
$ws.Cells.Item(7,17).Value2 = 0.05569913509533335
$ws.Cells.Item(8,17).Value2 = 0.1885220495682222
$ws.Cells.Item(9,17).Value2 = 0.01945336730866667
$ws.Cells.Item(10,17).Value2 = 0.1485197354044445
$ws.Cells.Item(11,17).Value2 = 0.006231808089111111

# Column R
$ws.Cells.Item(2,18).Value2 = 0.116559387552
$ws.Cells.Item(3,18).Value2 = 0.394512672416
$ws.Cells.Item(4,18).Value2 = 0.040709296032
$ws.Cells.Item(5,18).Value2 = 0.31080140416
$ws.Cells.Item(6,18).Value2 = 0.013041059488
$ws.Cells.Item(7,18).Value2 = 0.5012922158580001
$ws.Cells.Item(8,18).Value2 = 1.696698446114
$ws.Cells.Item(9,18).Value2 = 0.175080305778
$ws.Cells.Item(10,18).Value2 = 1.33667761864
$ws.Cells.Item(11,18).Value2 = 0.056086272802

# Column S
$ws.Cells.Item(2,19).Value2 = 0.02511266298316604
$ws.Cells.Item(3,19).Value2 = 0.0849975621273003
$ws.Cells.Item(4,19).Value2 = 0.008770797899718484
$ws.Cells.Item(5,19).Value2 = 0.06696201036474075
$ws.Cells.Item(6,19).Value2 = 0.002809689882073719
$ws.Cells.Item(7,19).Value2 = 0.1080031624849634
$ws.Cells.Item(8,19).Value2 = 0.3655528495490218
$ws.Cells.Item(9,19).Value2 = 0.03772096616440337
$ws.Cells.Item(10,19).Value2 = 0.2879865385280268
$ws.Cells.Item(11,19).Value2 = 0.01208376001658538

# Column T
$ws.Cells.Item(2,20).Value2 = 0.02511266298316604
$ws.Cells.Item(3,20).Value2 = 0.08499756212730028
$ws.Cells.Item(4,20).Value2 = 0.008770797899718482
$ws.Cells.Item(5,20).Value2 = 0.06696201036474075
$ws.Cells.Item(6,20).Value2 = 0.002809689882073719
$ws.Cells.Item(7,20).Value2 = 0.1080031624849634
$ws.Cells.Item(8,20).Value2 = 0.3655528495490217
$ws.Cells.Item(9,20).Value2 = 0.03772096616440336
$ws.Cells.Item(10,20).Value2 = 0.2879865385280268
$ws.Cells.Item(11,20).Value2 = 0.01208376001658538

